$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.802.05"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.750.04"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'236.05"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5081"
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D8").Value = "'40.78"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "'0.2677"
$ws.Range("E9").Value = "  +7.46%  "
$ws.Range("D10").Value = "'0.06190"
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").Value = "1.753.48"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'0.06935"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "'15.42"
$ws.Range("E13").Value = "  +4.12%  "
$ws.Range("D14").Value = "'0.6271"
$ws.Range("E14").Value = "  +10.81%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'77.65"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "25.806.78"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D21").Value = "'0.000006678"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "1.978.34"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'4.054"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'8.251"
$ws.Range("E24").Value = "  +5.00%  "
$ws.Range("D25").Value = "'5.128"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").Value = "'136.72"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'1.455"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'15.13"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("D29").Value = "'1.742"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").Value = "'102.69"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'0.08194"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'3.695"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'3.402"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "'0.04421"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.650"
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9968"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5987"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.633"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01565"
$ws.Range("E40").Value = "  +4.72%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.944"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'101.35"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.7503"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3817"
$ws.Range("E45").Value = "  +2.93%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'4.890"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05503"
$ws.Range("E47").Value = "  +7.65%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1097"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'5.938"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'30.05"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'52.72"
$ws.Range("E51").Value = "  +0.51%  "
